# Update the lattice-multiplication exercise table: every cell's problem
# (top line "A x B"), the two spaced-out digit factors, and the two partial
# seed digits down the left edge are replaced with a new set of values.
# The "  ----" divider line is identical before and after, so it is simply
# re-supplied unchanged as part of each cell's full text.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$NL = [char]11   # manual line break (<w:br/>) inside a Word Range.Text

# row, col, top line, factors line, dash line, and the two left-edge seed lines
$cells = @(
    @(1, 1, "84 x 43", "  4    3", "8|    |", "4|    |"),
    @(1, 2, "50 x 51", "  5    1", "5|    |", "0|    |"),
    @(1, 3, "18 x 69", "  6    9", "1|    |", "8|    |"),

    @(2, 1, "59 x 36", "  3    6", "5|    |", "9|    |"),
    @(2, 2, "25 x 90", "  9    0", "2|    |", "5|    |"),
    @(2, 3, "96 x 44", "  4    4", "9|    |", "6|    |"),

    @(3, 1, "84 x 44", "  4    4", "8|    |", "4|    |"),
    @(3, 2, "45 x 47", "  4    7", "4|    |", "5|    |"),
    @(3, 3, "89 x 45", "  4    5", "8|    |", "9|    |"),

    @(4, 1, "24 x 14", "  1    4", "2|    |", "4|    |"),
    @(4, 2, "13 x 29", "  2    9", "1|    |", "3|    |"),
    @(4, 3, "35 x 73", "  7    3", "3|    |", "5|    |"),

    @(5, 1, "67 x 29", "  2    9", "6|    |", "7|    |"),
    @(5, 2, "11 x 16", "  1    6", "1|    |", "1|    |"),
    @(5, 3, "45 x 63", "  6    3", "4|    |", "5|    |")
)

foreach ($row in $cells) {
    $r = [int]$row[0]
    $c = [int]$row[1]
    $top = $row[2]
    $factors = $row[3]
    $seedA = $row[4]
    $seedB = $row[5]

    $cell = $tbl.Cell($r, $c)
    $newText = $top + $NL + $factors + $NL + "  ----" + $NL + $seedA + $NL + $seedB
    $cell.Range.Text = $newText
}

Write-Output "updated $($cells.Length) cells"
